$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "copy xlsx config": duplicate the KIJIJIAUTOS sheet to create the new
#    FACEBOOK sheet, placed right after it (i.e. as the new last tab).
# ---------------------------------------------------------------------------
$src = $wb.Worksheets.Item("KIJIJIAUTOS")
$src.Copy($null, $src)
$new = $wb.Worksheets.Item($wb.Worksheets.Count)
$new.Name = "FACEBOOK"
$new.Activate()
$new.Range("A1").Select()

# KIJIJIAUTOS itself also ends up with a slightly different remembered
# selection after the copy operation.
$src.Activate()
$src.Range("G16").Select()

# ---------------------------------------------------------------------------
# 2. Clean up the mojibake / curly quotes around 'particuliers' in the
#    VO_CA expectation config (cell C47), replacing the two typographic
#    apostrophes (\u2019) with regular straight apostrophes.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("VO_CA")
$cell = $ws.Range("C47")
$cell.Value = "{'value_set' : ['pro','private','dealer','professional','particuliers'], 'mostly' : 0.95}"

# Re-apply explicit character-level formatting on the two replaced quote
# characters, matching how the retyped text kept its own font run.
$cell.Characters(57, 1).Font.Name = "Calibri"
$cell.Characters(70, 1).Font.Name = "Calibri"

# Retyping the cell also nudged the row height up slightly.
$ws.Rows.Item(47).RowHeight = 15

# ---------------------------------------------------------------------------
# 3. Restore VO_CA as the active sheet/tab with its updated last-selected
#    cell.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("C48").Select()
